# Refresh the cryptocurrency Price (column D) and Volume(1h) change
# (column E) figures with the latest scraped snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "27.387.72",
# "224.95"); format as Text first so Excel keeps these as literal
# strings instead of auto-converting them to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.387.72"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.712.53"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D5").Value = "224.95"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "0.5304"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.2669"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "0.06679"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").Value = "0.07699"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "4.518"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "1.948.40"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "1.709.07"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "0.5835"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "0.0₅8214"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "68.10"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "27.388.87"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "221.36"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "4.646"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "6.004"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("D26").Value = "1.687"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "0.1212"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "7.270"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").Value = "16.26"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "0.05361"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "1.299"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "3.459"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "1.644"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "0.9532"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "2.401"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "0.5869"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "0.01641"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "1.090.70"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").Value = "5.820"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").Value = "0.8472"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "101.16"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "1.854.76"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").Value = "57.97"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "0.4534"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "1.007"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "8.101"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "0.05241"
$ws.Range("E51").Value = "  -0.18%  "

# Restore the original (default/General) cell style on column D
# now that the text values are safely stored.
$priceRange.Style = "Normal"
